$d = $word.ActiveDocument

$replacements = @(
    @("939÷6=", "918÷2="),
    @("219÷6=", "990÷3="),
    @("871÷5=", "922÷3="),
    @("895÷8=", "835÷2="),
    @("881÷5=", "170÷7="),
    @("235÷6=", "980÷7="),
    @("954÷8=", "561÷5="),
    @("546÷2=", "402÷6="),
    @("266÷6=", "391÷4="),
    @("400÷7=", "501÷9="),
    @("573÷6=", "584÷6="),
    @("940÷5=", "319÷5="),
    @("975÷5=", "450÷2="),
    @("277÷6=", "823÷9="),
    @("292÷9=", "701÷2="),
    @("514÷5=", "134÷5="),
    @("867÷5=", "631÷8="),
    @("325÷6=", "111÷8="),
    @("641÷4=", "962÷4="),
    @("219÷7=", "151÷7="),
    @("232÷4=", "710÷4="),
    @("593÷7=", "492÷2="),
    @("815÷5=", "550÷7="),
    @("921÷6=", "326÷7="),
    @("217÷2=", "790÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
